$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The old lone "[]" string (previously at B38) is being retired - clear that
# cell first so the shared-string slot becomes free and can be recycled for
# the new text that replaces it.
$ws.Range("B38").ClearContents()

# --- "Verschiebung der Spindrehkurven" block -------------------------------
$ws.Range("A40").Value = "Verschiebung der Spindrehkurven"

$ws.Range("A41").Value = "max1"
$ws.Range("B41").Value = 0.0348261

$ws.Range("A42").Value = "max2"
$ws.Range("B42").Value = 0.0455357

$ws.Range("B43").Formula = "=(B42-B41)/2"

# --- "kl-fit" block ---------------------------------------------------------
$ws.Range("A45").Value = "kl-fit"

$ws.Range("A46").Value = "k="
$ws.Range("B46").Value = -0.064

$ws.Range("A47").Value = "l="
$ws.Range("B47").Value = -1.58

# --- "Einstrahlspule" block --------------------------------------------------
$ws.Range("A50").Value = "Einstrahlspule"
$ws.Range("B50").Value = "2,6 cm"
$ws.Range("B51").Value = "1,7cm"
$ws.Range("C51").Value = "Durchmesser"
$ws.Range("C50").Value = "Länge"

$ws.Range("A53").Value = "Halbwertsbreite"

# --- view state --------------------------------------------------------------
$ws.Range("A54").Select()
